$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '62.994.44'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '3.038.07'
$ws.Range('E3').Value = '  -3.43%  '
$ws.Range('E4').Value = '  +0.16%  '
Set-TextCell 'D5' '583.72'
$ws.Range('E5').Value = '  -1.28%  '
Set-TextCell 'D6' '153.19'
$ws.Range('E6').Value = '  +4.77%  '
$ws.Range('E7').Value = '  +0.05%  '
Set-TextCell 'D8' '0.535'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '3.031.60'
$ws.Range('E9').Value = '  -3.32%  '
Set-TextCell 'D10' '0.154'
$ws.Range('E10').Value = '  -5.28%  '
Set-TextCell 'D11' '5.77'
$ws.Range('E11').Value = '  -2.28%  '
Set-TextCell 'D12' '0.447'
$ws.Range('E12').Value = '  -2.19%  '
Set-TextCell 'D13' '36.50'
$ws.Range('E13').Value = '  -1.85%  '
Set-TextCell 'D14' '0.0000235'
$ws.Range('E14').Value = '  -5.17%  '
$ws.Range('D15').Value = '3.580.35'
$ws.Range('E15').Value = '  -2.33%  '
Set-TextCell 'D16' '0.118'
$ws.Range('E16').Value = '  -2.19%  '
$ws.Range('D17').Value = '63.281.33'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.069.95'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D19' '7.04'
$ws.Range('E19').Value = '  -3.00%  '
Set-TextCell 'D20' '470.21'
$ws.Range('E20').Value = '  +0.42%  '
Set-TextCell 'D21' '14.24'
$ws.Range('E21').Value = '  -0.82%  '
Set-TextCell 'D22' '0.699'
$ws.Range('E22').Value = '  -4.43%  '
Set-TextCell 'D23' '7.40'
$ws.Range('E23').Value = '  -2.34%  '
Set-TextCell 'D24' '2.37'
$ws.Range('E24').Value = '  -0.55%  '
Set-TextCell 'D25' '80.23'
$ws.Range('E25').Value = '  -0.86%  '
Set-TextCell 'D26' '12.66'
$ws.Range('E26').Value = '  -3.74%  '
Set-TextCell 'D27' '10.18'
$ws.Range('E27').Value = '  +3.64%  '
Set-TextCell 'D28' '0.996'
$ws.Range('E28').Value = '  -0.44%  '
Set-TextCell 'D29' '7.45'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 'D30' '1.01'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D31' '2.63'
$ws.Range('E31').Value = '  -2.84%  '
Set-TextCell 'D32' '2.13'
$ws.Range('E32').Value = '  -3.55%  '
$ws.Range('E33').Value = '  -2.12%  '
Set-TextCell 'D34' '26.82'
$ws.Range('E34').Value = '  -2.87%  '
$ws.Range('D35').Value = '0.0₃0812'
$ws.Range('E35').Value = '  -5.76%  '
Set-TextCell 'D36' '1.04'
$ws.Range('E36').Value = '  -2.36%  '
Set-TextCell 'D37' '3.32'
$ws.Range('E37').Value = '  +1.72%  '
Set-TextCell 'D38' '5.92'
$ws.Range('E38').Value = '  -3.86%  '
Set-TextCell 'D39' '2.17'
$ws.Range('E39').Value = '  -4.22%  '
Set-TextCell 'D40' '9.23'
$ws.Range('E40').Value = '  -1.52%  '
Set-TextCell 'D41' '50.40'
$ws.Range('E41').Value = '  -1.83%  '
Set-TextCell 'D42' '434.10'
$ws.Range('E42').Value = '  -6.08%  '
Set-TextCell 'D43' '0.283'
$ws.Range('E43').Value = '  -3.15%  '
Set-TextCell 'D44' '40.91'
$ws.Range('E44').Value = '  +2.34%  '
Set-TextCell 'D45' '0.111'
$ws.Range('E45').Value = '  +3.60%  '
Set-TextCell 'D46' '0.0355'
$ws.Range('E46').Value = '  -4.73%  '
$ws.Range('D47').Value = '2.776.88'
$ws.Range('E47').Value = '  -3.98%  '
Set-TextCell 'D48' '130.01'
$ws.Range('E48').Value = '  -2.12%  '
Set-TextCell 'D49' '0.999'
$ws.Range('E49').Value = '  +0.06%  '
Set-TextCell 'D50' '24.81'
$ws.Range('E50').Value = '  +3.30%  '
Set-TextCell 'D51' '2.21'
$ws.Range('E51').Value = '  -0.89%  '
